$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix out-of-order rows caused by inserting the new match below (rows shifted) ---

# Row 143
$ws.Range("B143").Value = "england"
$ws.Range("C143").Value = "championship"
$ws.Range("D143").Value = "2023-2024"
$ws.Range("E143").Value = 45220.66666666666
$ws.Range("F143").Value = "Blackburn"
$ws.Range("G143").Value = 1
$ws.Range("H143").Value = "Cardiff"
$ws.Range("I143").Value = 0
$ws.Range("J143").Value = 2.01
$ws.Range("K143").Value = "13/10/2023 15:43"
$ws.Range("L143").Value = 2.05
$ws.Range("M143").Value = "21/10/2023 15:58"
$ws.Range("N143").Value = 3.78
$ws.Range("O143").Value = "13/10/2023 15:43"
$ws.Range("P143").Value = 3.72
$ws.Range("Q143").Value = "21/10/2023 15:38"
$ws.Range("R143").Value = 3.67
$ws.Range("S143").Value = "13/10/2023 15:43"
$ws.Range("T143").Value = 3.68
$ws.Range("U143").Value = "21/10/2023 15:58"
$ws.Range("V143").Value = "https://www.betexplorer.com/football/england/championship/blackburn-cardiff/KMhLBBp8/"

# Row 144
$ws.Range("B144").Value = "england"
$ws.Range("C144").Value = "championship"
$ws.Range("D144").Value = "2023-2024"
$ws.Range("E144").Value = 45220.66666666666
$ws.Range("F144").Value = "Hull"
$ws.Range("G144").Value = 1
$ws.Range("H144").Value = "Southampton"
$ws.Range("I144").Value = 2
$ws.Range("J144").Value = 2.71
$ws.Range("K144").Value = "09/10/2023 14:42"
$ws.Range("L144").Value = 3.06
$ws.Range("M144").Value = "21/10/2023 15:52"
$ws.Range("N144").Value = 3.53
$ws.Range("O144").Value = "09/10/2023 14:42"
$ws.Range("P144").Value = 3.55
$ws.Range("Q144").Value = "21/10/2023 15:52"
$ws.Range("R144").Value = 2.62
$ws.Range("S144").Value = "09/10/2023 14:42"
$ws.Range("T144").Value = 2.38
$ws.Range("U144").Value = "21/10/2023 15:59"
$ws.Range("V144").Value = "https://www.betexplorer.com/football/england/championship/hull-city-southampton/Q1ddTLgK/"

# Row 159
$ws.Range("B159").Value = "england"
$ws.Range("C159").Value = "championship"
$ws.Range("D159").Value = "2023-2024"
$ws.Range("E159").Value = 45227.66666666666
$ws.Range("F159").Value = "Watford"
$ws.Range("G159").Value = 2
$ws.Range("H159").Value = "Millwall"
$ws.Range("I159").Value = 2
$ws.Range("J159").Value = 2.07
$ws.Range("K159").Value = "21/10/2023 18:13"
$ws.Range("L159").Value = 2.22
$ws.Range("M159").Value = "28/10/2023 15:56"
$ws.Range("N159").Value = 3.48
$ws.Range("O159").Value = "21/10/2023 18:13"
$ws.Range("P159").Value = 3.22
$ws.Range("Q159").Value = "28/10/2023 15:56"
$ws.Range("R159").Value = 3.79
$ws.Range("S159").Value = "21/10/2023 18:13"
$ws.Range("T159").Value = 3.74
$ws.Range("U159").Value = "28/10/2023 15:56"
$ws.Range("V159").Value = "https://www.betexplorer.com/football/england/championship/watford-millwall/fuq5rhaO/"

# Row 162
$ws.Range("B162").Value = "england"
$ws.Range("C162").Value = "championship"
$ws.Range("D162").Value = "2023-2024"
$ws.Range("E162").Value = 45227.66666666666
$ws.Range("F162").Value = "Middlesbrough"
$ws.Range("G162").Value = 0
$ws.Range("H162").Value = "Stoke"
$ws.Range("I162").Value = 2
$ws.Range("J162").Value = 1.88
$ws.Range("K162").Value = "21/10/2023 18:13"
$ws.Range("L162").Value = 1.78
$ws.Range("M162").Value = "28/10/2023 15:51"
$ws.Range("N162").Value = 3.84
$ws.Range("O162").Value = "21/10/2023 18:13"
$ws.Range("P162").Value = 3.91
$ws.Range("Q162").Value = "28/10/2023 15:56"
$ws.Range("R162").Value = 4.16
$ws.Range("S162").Value = "21/10/2023 18:13"
$ws.Range("T162").Value = 4.69
$ws.Range("U162").Value = "28/10/2023 15:51"
$ws.Range("V162").Value = "https://www.betexplorer.com/football/england/championship/middlesbrough-stoke-city/MiGhmGpg/"

# Row 163
$ws.Range("B163").Value = "england"
$ws.Range("C163").Value = "championship"
$ws.Range("D163").Value = "2023-2024"
$ws.Range("E163").Value = 45227.66666666666
$ws.Range("F163").Value = "Hull"
$ws.Range("G163").Value = 1
$ws.Range("H163").Value = "Preston"
$ws.Range("I163").Value = 0
$ws.Range("J163").Value = 2.09
$ws.Range("K163").Value = "21/10/2023 18:13"
$ws.Range("L163").Value = 2.01
$ws.Range("M163").Value = "28/10/2023 15:51"
$ws.Range("N163").Value = 3.52
$ws.Range("O163").Value = "21/10/2023 18:13"
$ws.Range("P163").Value = 3.5
$ws.Range("Q163").Value = "28/10/2023 15:38"
$ws.Range("R163").Value = 3.73
$ws.Range("S163").Value = "21/10/2023 18:13"
$ws.Range("T163").Value = 4.07
$ws.Range("U163").Value = "28/10/2023 15:51"
$ws.Range("V163").Value = "https://www.betexplorer.com/football/england/championship/hull-city-preston/OGzIgI7P/"

# Row 164
$ws.Range("B164").Value = "england"
$ws.Range("C164").Value = "championship"
$ws.Range("D164").Value = "2023-2024"
$ws.Range("E164").Value = 45227.66666666666
$ws.Range("F164").Value = "Cardiff"
$ws.Range("G164").Value = 2
$ws.Range("H164").Value = "Bristol City"
$ws.Range("I164").Value = 0
$ws.Range("J164").Value = 2.38
$ws.Range("K164").Value = "21/10/2023 18:13"
$ws.Range("L164").Value = 2.04
$ws.Range("M164").Value = "28/10/2023 15:58"
$ws.Range("N164").Value = 3.39
$ws.Range("O164").Value = "21/10/2023 18:13"
$ws.Range("P164").Value = 3.68
$ws.Range("Q164").Value = "28/10/2023 15:58"
$ws.Range("R164").Value = 3.19
$ws.Range("S164").Value = "21/10/2023 18:13"
$ws.Range("T164").Value = 3.74
$ws.Range("U164").Value = "28/10/2023 15:58"
$ws.Range("V164").Value = "https://www.betexplorer.com/football/england/championship/cardiff-bristol-city/lbbd3mhl/"

# Row 170
$ws.Range("B170").Value = "england"
$ws.Range("C170").Value = "championship"
$ws.Range("D170").Value = "2023-2024"
$ws.Range("E170").Value = 45234.66666666666
$ws.Range("F170").Value = "West Brom"
$ws.Range("G170").Value = 3
$ws.Range("H170").Value = "Hull"
$ws.Range("I170").Value = 1
$ws.Range("J170").Value = 1.85
$ws.Range("K170").Value = "28/10/2023 18:13"
$ws.Range("L170").Value = 1.95
$ws.Range("M170").Value = "04/11/2023 15:59"
$ws.Range("N170").Value = 3.72
$ws.Range("O170").Value = "28/10/2023 18:13"
$ws.Range("P170").Value = 3.55
$ws.Range("Q170").Value = "04/11/2023 15:53"
$ws.Range("R170").Value = 4.31
$ws.Range("S170").Value = "28/10/2023 18:13"
$ws.Range("T170").Value = 4.27
$ws.Range("U170").Value = "04/11/2023 15:59"
$ws.Range("V170").Value = "https://www.betexplorer.com/football/england/championship/west-brom-hull-city/2qRCvqbj/"

# Row 171
$ws.Range("B171").Value = "england"
$ws.Range("C171").Value = "championship"
$ws.Range("D171").Value = "2023-2024"
$ws.Range("E171").Value = 45234.66666666666
$ws.Range("F171").Value = "Swansea"
$ws.Range("G171").Value = 0
$ws.Range("H171").Value = "Sunderland"
$ws.Range("I171").Value = 0
$ws.Range("J171").Value = 2.6
$ws.Range("K171").Value = "28/10/2023 18:13"
$ws.Range("L171").Value = 3.03
$ws.Range("M171").Value = "04/11/2023 15:31"
$ws.Range("N171").Value = 3.46
$ws.Range("O171").Value = "28/10/2023 18:13"
$ws.Range("P171").Value = 3.5
$ws.Range("Q171").Value = "04/11/2023 15:30"
$ws.Range("R171").Value = 2.74
$ws.Range("S171").Value = "28/10/2023 18:13"
$ws.Range("T171").Value = 2.43
$ws.Range("U171").Value = "04/11/2023 15:31"
$ws.Range("V171").Value = "https://www.betexplorer.com/football/england/championship/swansea-sunderland/ADn51RM0/"

# Row 172
$ws.Range("B172").Value = "england"
$ws.Range("C172").Value = "championship"
$ws.Range("D172").Value = "2023-2024"
$ws.Range("E172").Value = 45234.66666666666
$ws.Range("F172").Value = "Stoke"
$ws.Range("G172").Value = 0
$ws.Range("H172").Value = "Cardiff"
$ws.Range("I172").Value = 0
$ws.Range("J172").Value = 2.03
$ws.Range("K172").Value = "28/10/2023 18:13"
$ws.Range("L172").Value = 2.11
$ws.Range("M172").Value = "04/11/2023 15:12"
$ws.Range("N172").Value = 3.61
$ws.Range("O172").Value = "28/10/2023 18:13"
$ws.Range("P172").Value = 3.41
$ws.Range("Q172").Value = "04/11/2023 15:34"
$ws.Range("R172").Value = 3.7
$ws.Range("S172").Value = "28/10/2023 18:13"
$ws.Range("T172").Value = 3.84
$ws.Range("U172").Value = "04/11/2023 15:12"
$ws.Range("V172").Value = "https://www.betexplorer.com/football/england/championship/stoke-city-cardiff/WMm1277f/"

# Row 173
$ws.Range("B173").Value = "england"
$ws.Range("C173").Value = "championship"
$ws.Range("D173").Value = "2023-2024"
$ws.Range("E173").Value = 45234.66666666666
$ws.Range("F173").Value = "Rotherham"
$ws.Range("G173").Value = 1
$ws.Range("H173").Value = "QPR"
$ws.Range("I173").Value = 1
$ws.Range("J173").Value = 2.61
$ws.Range("K173").Value = "28/10/2023 18:13"
$ws.Range("L173").Value = 3.15
$ws.Range("M173").Value = "04/11/2023 15:56"
$ws.Range("N173").Value = 3.3
$ws.Range("O173").Value = "28/10/2023 18:13"
$ws.Range("P173").Value = 3.31
$ws.Range("Q173").Value = "04/11/2023 15:56"
$ws.Range("R173").Value = 2.9
$ws.Range("S173").Value = "28/10/2023 18:13"
$ws.Range("T173").Value = 2.45
$ws.Range("U173").Value = "04/11/2023 15:56"
$ws.Range("V173").Value = "https://www.betexplorer.com/football/england/championship/rotherham-qpr/p8W7uPrp/"

# Row 174
$ws.Range("B174").Value = "england"
$ws.Range("C174").Value = "championship"
$ws.Range("D174").Value = "2023-2024"
$ws.Range("E174").Value = 45234.66666666666
$ws.Range("F174").Value = "Preston"
$ws.Range("G174").Value = 3
$ws.Range("H174").Value = "Coventry"
$ws.Range("I174").Value = 2
$ws.Range("J174").Value = 2.7
$ws.Range("K174").Value = "28/10/2023 18:12"
$ws.Range("L174").Value = 2.86
$ws.Range("M174").Value = "04/11/2023 15:56"
$ws.Range("N174").Value = 3.31
$ws.Range("O174").Value = "28/10/2023 18:12"
$ws.Range("P174").Value = 3.33
$ws.Range("Q174").Value = "04/11/2023 15:56"
$ws.Range("R174").Value = 2.82
$ws.Range("S174").Value = "28/10/2023 18:12"
$ws.Range("T174").Value = 2.65
$ws.Range("U174").Value = "04/11/2023 15:56"
$ws.Range("V174").Value = "https://www.betexplorer.com/football/england/championship/preston-coventry/8r9vciEH/"

# Row 175
$ws.Range("B175").Value = "england"
$ws.Range("C175").Value = "championship"
$ws.Range("D175").Value = "2023-2024"
$ws.Range("E175").Value = 45234.66666666666
$ws.Range("F175").Value = "Huddersfield"
$ws.Range("G175").Value = 0
$ws.Range("H175").Value = "Watford"
$ws.Range("I175").Value = 0
$ws.Range("J175").Value = 2.83
$ws.Range("K175").Value = "28/10/2023 18:13"
$ws.Range("L175").Value = 3.67
$ws.Range("M175").Value = "04/11/2023 15:59"
$ws.Range("N175").Value = 3.4
$ws.Range("O175").Value = "28/10/2023 18:13"
$ws.Range("P175").Value = 3.47
$ws.Range("Q175").Value = "04/11/2023 15:59"
$ws.Range("R175").Value = 2.56
$ws.Range("S175").Value = "28/10/2023 18:13"
$ws.Range("T175").Value = 2.14
$ws.Range("U175").Value = "04/11/2023 15:59"
$ws.Range("V175").Value = "https://www.betexplorer.com/football/england/championship/huddersfield-watford/CdDJ1Fan/"

# Row 176
$ws.Range("B176").Value = "england"
$ws.Range("C176").Value = "championship"
$ws.Range("D176").Value = "2023-2024"
$ws.Range("E176").Value = 45234.66666666666
$ws.Range("F176").Value = "Millwall"
$ws.Range("G176").Value = 0
$ws.Range("H176").Value = "Southampton"
$ws.Range("I176").Value = 1
$ws.Range("J176").Value = 2.93
$ws.Range("K176").Value = "28/10/2023 18:13"
$ws.Range("L176").Value = 3.65
$ws.Range("M176").Value = "04/11/2023 15:56"
$ws.Range("N176").Value = 3.52
$ws.Range("O176").Value = "28/10/2023 18:13"
$ws.Range("P176").Value = 3.61
$ws.Range("Q176").Value = "04/11/2023 15:59"
$ws.Range("R176").Value = 2.42
$ws.Range("S176").Value = "28/10/2023 18:13"
$ws.Range("T176").Value = 2.09
$ws.Range("U176").Value = "04/11/2023 15:59"
$ws.Range("V176").Value = "https://www.betexplorer.com/football/england/championship/millwall-southampton/KMGRagUb/"

# Row 177
$ws.Range("B177").Value = "england"
$ws.Range("C177").Value = "championship"
$ws.Range("D177").Value = "2023-2024"
$ws.Range("E177").Value = 45234.66666666666
$ws.Range("F177").Value = "Bristol City"
$ws.Range("G177").Value = 1
$ws.Range("H177").Value = "Sheffield Wed"
$ws.Range("I177").Value = 0
$ws.Range("J177").Value = 1.86
$ws.Range("K177").Value = "28/10/2023 18:12"
$ws.Range("L177").Value = 2.24
$ws.Range("M177").Value = "04/11/2023 15:59"
$ws.Range("N177").Value = 3.8
$ws.Range("O177").Value = "28/10/2023 18:12"
$ws.Range("P177").Value = 3.21
$ws.Range("Q177").Value = "04/11/2023 15:59"
$ws.Range("R177").Value = 4.36
$ws.Range("S177").Value = "28/10/2023 18:12"
$ws.Range("T177").Value = 3.71
$ws.Range("U177").Value = "04/11/2023 15:59"
$ws.Range("V177").Value = "https://www.betexplorer.com/football/england/championship/bristol-city-sheffield-wed/Y5EF2ept/"

# Row 178
$ws.Range("B178").Value = "england"
$ws.Range("C178").Value = "championship"
$ws.Range("D178").Value = "2023-2024"
$ws.Range("E178").Value = 45234.66666666666
$ws.Range("F178").Value = "Birmingham"
$ws.Range("G178").Value = 2
$ws.Range("H178").Value = "Ipswich"
$ws.Range("I178").Value = 2
$ws.Range("J178").Value = 3.73
$ws.Range("K178").Value = "28/10/2023 18:13"
$ws.Range("L178").Value = 4.03
$ws.Range("M178").Value = "04/11/2023 15:57"
$ws.Range("N178").Value = 3.63
$ws.Range("O178").Value = "28/10/2023 18:13"
$ws.Range("P178").Value = 3.81
$ws.Range("Q178").Value = "04/11/2023 15:55"
$ws.Range("R178").Value = 2.01
$ws.Range("S178").Value = "28/10/2023 18:13"
$ws.Range("T178").Value = 1.93
$ws.Range("U178").Value = "04/11/2023 15:32"
$ws.Range("V178").Value = "https://www.betexplorer.com/football/england/championship/birmingham-ipswich/8rm9sCEU/"

# --- Append new match row 183 (Sunderland vs Birmingham) ---

# Copy formatting from the last existing data row (182) for the indexed/date columns
$ws.Range("A182").Copy()
$ws.Range("A183").PasteSpecial(-4122)
$ws.Range("E182").Copy()
$ws.Range("E183").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A183").Value = 182
$ws.Range("B183").Value = "england"
$ws.Range("C183").Value = "championship"
$ws.Range("D183").Value = "2023-2024"
$ws.Range("E183").Value = 45241.5625
$ws.Range("F183").Value = "Sunderland"
$ws.Range("G183").Value = 3
$ws.Range("H183").Value = "Birmingham"
$ws.Range("I183").Value = 1
$ws.Range("J183").Value = 1.71
$ws.Range("K183").Value = "04/11/2023 13:42"
$ws.Range("L183").Value = 1.65
$ws.Range("M183").Value = "11/11/2023 13:20"
$ws.Range("N183").Value = 3.95
$ws.Range("O183").Value = "04/11/2023 13:42"
$ws.Range("P183").Value = 4.2
$ws.Range("Q183").Value = "11/11/2023 13:20"
$ws.Range("R183").Value = 5.04
$ws.Range("S183").Value = "04/11/2023 13:42"
$ws.Range("T183").Value = 5.32
$ws.Range("U183").Value = "11/11/2023 13:25"
$ws.Range("V183").Value = "https://www.betexplorer.com/football/england/championship/sunderland-birmingham/IwCrFpi0/"
